# Updates the "cryptos" price/volume snapshot (GitHub Actions refresh).
# Values in columns D/E (and some in B/C) are stored as *text*, not numbers
# (several, e.g. "1.002", "30.506.18", look numeric but use "." as a
# thousands separator / contain multiple dots, so they must stay strings).
# Forcing NumberFormat to "@" (Text) before assigning the value keeps Excel
# from re-interpreting numeric-looking strings as numbers; resetting the
# Style back to "Normal" afterwards avoids leaving a stray text-format style
# applied to the cell, matching the unstyled cells in the original sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.506.18"
Set-TextValue $ws.Range("E2") "  +0.58%  "

Set-TextValue $ws.Range("D3") "2.108.32"
Set-TextValue $ws.Range("E3") "  +4.94%  "

Set-TextValue $ws.Range("E4") "  +0.01%  "

Set-TextValue $ws.Range("D5") "329.28"
Set-TextValue $ws.Range("E5") "  +1.39%  "

Set-TextValue $ws.Range("D6") "1.002"
Set-TextValue $ws.Range("E6") "  +0.09%  "

Set-TextValue $ws.Range("D7") "0.5264"
Set-TextValue $ws.Range("E7") "  +2.64%  "

Set-TextValue $ws.Range("D8") "0.4370"
Set-TextValue $ws.Range("E8") "  +2.57%  "

Set-TextValue $ws.Range("D9") "0.08882"
Set-TextValue $ws.Range("E9") "  +2.04%  "

Set-TextValue $ws.Range("D10") "47.28"
Set-TextValue $ws.Range("E10") "  +9.56%  "

Set-TextValue $ws.Range("E11") "  +2.43%  "

Set-TextValue $ws.Range("D12") "24.54"
Set-TextValue $ws.Range("E12") "  -0.74%  "

Set-TextValue $ws.Range("D13") "2.103.69"
Set-TextValue $ws.Range("E13") "  +4.73%  "

Set-TextValue $ws.Range("D14") "6.723"
Set-TextValue $ws.Range("E14") "  +2.29%  "

Set-TextValue $ws.Range("D15") "7.761"
Set-TextValue $ws.Range("E15") "  +3.85%  "

Set-TextValue $ws.Range("D16") "96.37"
Set-TextValue $ws.Range("E16") "  +2.13%  "

Set-TextValue $ws.Range("D17") "1.003"
Set-TextValue $ws.Range("E17") "  +0.11%  "

Set-TextValue $ws.Range("E18") "  +1.18%  "

Set-TextValue $ws.Range("D19") "0.06645"
Set-TextValue $ws.Range("E19") "  +1.57%  "

Set-TextValue $ws.Range("E20") "  +0.96%  "

Set-TextValue $ws.Range("E21") "  +0.10%  "

Set-TextValue $ws.Range("D22") "6.331"
Set-TextValue $ws.Range("E22") "  +2.12%  "

Set-TextValue $ws.Range("D23") "30.549.01"
Set-TextValue $ws.Range("E23") "  +0.54%  "

Set-TextValue $ws.Range("D24") "12.30"
Set-TextValue $ws.Range("E24") "  +4.11%  "

Set-TextValue $ws.Range("D25") "2.355"
Set-TextValue $ws.Range("E25") "  +4.28%  "

Set-TextValue $ws.Range("D26") "2.354.13"
Set-TextValue $ws.Range("E26") "  +4.90%  "

Set-TextValue $ws.Range("D27") "22.44"
Set-TextValue $ws.Range("E27") "  +0.19%  "

Set-TextValue $ws.Range("E28") "  +7.03%  "

Set-TextValue $ws.Range("D29") "162.15"
Set-TextValue $ws.Range("E29") "  -0.07%  "

Set-TextValue $ws.Range("D30") "132.68"
Set-TextValue $ws.Range("E30") "  +1.23%  "

Set-TextValue $ws.Range("D31") "1.207"
Set-TextValue $ws.Range("E31") "  +6.18%  "

Set-TextValue $ws.Range("E32") "  +2.33%  "

Set-TextValue $ws.Range("D33") "1.667"
Set-TextValue $ws.Range("E33") "  +21.85%  "

Set-TextValue $ws.Range("E34") "  +2.11%  "

Set-TextValue $ws.Range("D35") "3.922"
Set-TextValue $ws.Range("E35") "  +2.53%  "

Set-TextValue $ws.Range("E36") "  +9.99%  "

Set-TextValue $ws.Range("D37") "0.02585"
Set-TextValue $ws.Range("E37") "  +2.43%  "

Set-TextValue $ws.Range("D38") "5.492"
Set-TextValue $ws.Range("E38") "  +0.58%  "

Set-TextValue $ws.Range("B39") "Aptos"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D39") "12.72"
Set-TextValue $ws.Range("E39") "  +2.91%  "

Set-TextValue $ws.Range("B40") "Hedera"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D40") "0.06679"
Set-TextValue $ws.Range("E40") "  -0.04%  "

Set-TextValue $ws.Range("D41") "0.2284"
Set-TextValue $ws.Range("E41") "  +4.21%  "

Set-TextValue $ws.Range("D42") "0.6801"
Set-TextValue $ws.Range("E42") "  +2.47%  "

Set-TextValue $ws.Range("D43") "1.273"
Set-TextValue $ws.Range("E43") "  +2.75%  "

Set-TextValue $ws.Range("E44") "  +0.07%  "

Set-TextValue $ws.Range("D45") "14.03"
Set-TextValue $ws.Range("E45") "  +3.25%  "

Set-TextValue $ws.Range("D46") "0.6386"
Set-TextValue $ws.Range("E46") "  +3.56%  "

Set-TextValue $ws.Range("D47") "2.211"
Set-TextValue $ws.Range("E47") "  +1.28%  "

Set-TextValue $ws.Range("D48") "3.622"
Set-TextValue $ws.Range("E48") "  -1.07%  "

Set-TextValue $ws.Range("D49") "1.251"
Set-TextValue $ws.Range("E49") "  -0.77%  "

Set-TextValue $ws.Range("B50") "WEMIXTOKEN"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D50") "1.198"
Set-TextValue $ws.Range("E50") "  +8.17%  "

Set-TextValue $ws.Range("B51") "Aave"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D51") "82.60"
Set-TextValue $ws.Range("E51") "  +2.35%  "
